$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price column (D) values
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.083.90"
$ws.Range("D2").ClearFormats()
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.778.29"
$ws.Range("D3").ClearFormats()
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.59"
$ws.Range("D5").ClearFormats()
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.778.70"
$ws.Range("D7").ClearFormats()
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.60"
$ws.Range("D14").ClearFormats()
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.410.42"
$ws.Range("D15").ClearFormats()
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.775.89"
$ws.Range("D16").ClearFormats()
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.71"
$ws.Range("D17").ClearFormats()
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "68.037.15"
$ws.Range("D18").ClearFormats()
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.55"
$ws.Range("D21").ClearFormats()
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.07"
$ws.Range("D22").ClearFormats()
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "0.0000148"
$ws.Range("D24").ClearFormats()
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.23"
$ws.Range("D26").ClearFormats()
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.33"
$ws.Range("D28").ClearFormats()
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.925.73"
$ws.Range("D31").ClearFormats()
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.56"
$ws.Range("D32").ClearFormats()
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "30.42"
$ws.Range("D33").ClearFormats()
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.734.99"
$ws.Range("D36").ClearFormats()
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.77"
$ws.Range("D37").ClearFormats()
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.138"
$ws.Range("D39").ClearFormats()
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.00"
$ws.Range("D42").ClearFormats()
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "406.71"
$ws.Range("D47").ClearFormats()
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "45.57"
$ws.Range("D48").ClearFormats()
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "143.52"
$ws.Range("D49").ClearFormats()
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.000272"
$ws.Range("D50").ClearFormats()
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "39.92"
$ws.Range("D51").ClearFormats()

# Update Volume(1h) column (E) values
$ws.Range("E2").Value = "  -0.86%  "
$ws.Range("E3").Value = "  -2.29%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("E5").Value = "  -0.95%  "
$ws.Range("E6").Value = "  -1.71%  "
$ws.Range("E7").Value = "  -2.33%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("E9").Value = "  -0.70%  "
$ws.Range("E10").Value = "  -2.86%  "
$ws.Range("E12").Value = "  -2.38%  "
$ws.Range("E13").Value = "  -2.36%  "
$ws.Range("E14").Value = "  -1.51%  "
$ws.Range("E15").Value = "  -2.32%  "
$ws.Range("E16").Value = "  -3.13%  "
$ws.Range("E17").Value = "  +2.04%  "
$ws.Range("E18").Value = "  -0.99%  "
$ws.Range("E19").Value = "  -3.36%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  -4.38%  "
$ws.Range("E22").Value = "  -1.12%  "
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("E24").Value = "  -8.64%  "
$ws.Range("E25").Value = "  -0.08%  "
$ws.Range("E26").Value = "  -1.59%  "
$ws.Range("E28").Value = "  -1.62%  "
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("E30").Value = "  -1.21%  "
$ws.Range("E31").Value = "  -2.29%  "
$ws.Range("E32").Value = "  -2.92%  "
$ws.Range("E33").Value = "  -2.97%  "
$ws.Range("E34").Value = "  -4.34%  "
$ws.Range("E35").Value = "  -2.24%  "
$ws.Range("E37").Value = "  -3.97%  "
$ws.Range("E38").Value = "  -1.78%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("E40").Value = "  -1.40%  "
$ws.Range("E41").Value = "  -2.73%  "
$ws.Range("E42").Value = "  -0.07%  "
$ws.Range("E43").Value = "  -2.87%  "
$ws.Range("E45").Value = "  -1.61%  "
$ws.Range("E47").Value = "  -3.17%  "
$ws.Range("E48").Value = "  -2.31%  "
$ws.Range("E49").Value = "  +0.55%  "
$ws.Range("E50").Value = "  -9.84%  "
$ws.Range("E51").Value = "  +3.64%  "
